{"js": "// Apply the \"Basic Weapons And Enemies\" content update:\n//   - Rename several enemy / weapon list items (prefix each with a new\n//     creature / item name while keeping the original parenthetical\n//     description).\n//   - Add two new enemy entries (Pyromaniac, Necromancer).\n//   - Add one new weapon entry (Sword).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map of the *original* full paragraph text -> new full paragraph text.\n// Using the original text as the lookup key keeps this robust against\n// paragraph re-indexing and avoids hard-coding paragraph positions.\nconst renames = new Map([\n  [\n    \"Standard Enemy (Moves towards player, uses standard weapon, low health)\",\n    \"Standard Skeleton (Moves towards player, uses standard weapon, low health)\",\n  ],\n  [\n    \"Elite Enemy (Moves towards player, uses Burst-Fire weapon, average health)\",\n    \"Mage (Moves towards player, uses Burst-Fire weapon, average health)\",\n  ],\n  [\n    \"Scaredy-Cat Enemy (Moves away from player, uses Rapid-Fire weapon, low health)\",\n    \"Ghost (Moves away from player, uses Rapid-Fire weapon, low health)\",\n  ],\n  [\n    \"Tank Enemy (Moves towards player, uses Spread Weapon, high health but slow movement)\",\n    \"Staff Ogre (Moves towards player, uses Staff, high health but slow movement)\",\n  ],\n  [\n    \"Standard Weapon (1 projectile per shot, average attack speed, average damage)\",\n    \"Wand (Standard Weapon) (1 projectile per shot, average attack speed, average damage)\",\n  ],\n  [\n    \"Spread weapon (Multiple projectiles per shot (3-5?), slow attack speed, average damage)\",\n    \"Staff (Multiple projectiles per shot (3-5?), slow attack speed, average damage)\",\n  ],\n  [\n    \"Explosive launcher (1 projectile per shot (explodes), (very?) slow attack speed, high damage to enemies hit directly, lower damage to enemies further away from projectile when it collides with something)\",\n    \"Fireball (1 projectile per shot (explodes), (very?) slow attack speed, high damage to enemies hit directly, lower damage to enemies further away from projectile when it collides with something)\",\n  ],\n]);\n\nlet tankOgreParagraph = null;\nlet burstFireParagraph = null;\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const para = items[i];\n  const text = para.text;\n  if (renames.has(text)) {\n    const newText = renames.get(text);\n    para.insertText(newText, \"Replace\");\n    if (text.indexOf(\"Tank Enemy\") === 0) {\n      tankOgreParagraph = para;\n    }\n  } else if (text.indexOf(\"Burst-Fire (1 projectile per shot\") === 0) {\n    burstFireParagraph = para;\n  }\n}\nawait context.sync();\n\n// New enemy entries, inserted right after the (renamed) \"Staff Ogre\" entry.\nif (tankOgreParagraph) {\n  const pyromaniac = tankOgreParagraph.insertParagraph(\n    \"Pyromaniac (Moves towards player quickly, uses Flamethrower, average health and fast movement)\",\n    \"After\"\n  );\n  pyromaniac.insertParagraph(\n    \"Necromancer (Moves away from player slowly, spawns Standard Skeletons and shoots a ring of projectiles around itself occasionally, High health and slow movement)\",\n    \"After\"\n  );\n}\n\n// New weapon entry, inserted right after the \"Burst-Fire\" entry.\nif (burstFireParagraph) {\n  burstFireParagraph.insertParagraph(\n    \"Sword (Melee weapon, hits enemies directly in front of player, fast attack speed, high damage)\",\n    \"After\"\n  );\n}\n\nawait context.sync();\n", "ps1": "# Apply the \"Basic Weapons And Enemies\" content update:\n#   - Rename several enemy / weapon list items (prefix each with a new\n#     creature / item name while keeping the original parenthetical\n#     description).\n#   - Add two new enemy entries (Pyromaniac, Necromancer).\n#   - Add one new weapon entry (Sword).\n$d = $word.ActiveDocument\n\n# Map of the *original* full paragraph text -> new full paragraph text.\n# Matching on the original text (rather than a hard-coded paragraph index)\n# keeps this robust to any incidental paragraph re-numbering.\n$renames = @{\n    \"Standard Enemy (Moves towards player, uses standard weapon, low health)\" = \"Standard Skeleton (Moves towards player, uses standard weapon, low health)\";\n    \"Elite Enemy (Moves towards player, uses Burst-Fire weapon, average health)\" = \"Mage (Moves towards player, uses Burst-Fire weapon, average health)\";\n    \"Scaredy-Cat Enemy (Moves away from player, uses Rapid-Fire weapon, low health)\" = \"Ghost (Moves away from player, uses Rapid-Fire weapon, low health)\";\n    \"Tank Enemy (Moves towards player, uses Spread Weapon, high health but slow movement)\" = \"Staff Ogre (Moves towards player, uses Staff, high health but slow movement)\";\n    \"Standard Weapon (1 projectile per shot, average attack speed, average damage)\" = \"Wand (Standard Weapon) (1 projectile per shot, average attack speed, average damage)\";\n    \"Spread weapon (Multiple projectiles per shot (3-5?), slow attack speed, average damage)\" = \"Staff (Multiple projectiles per shot (3-5?), slow attack speed, average damage)\";\n    \"Explosive launcher (1 projectile per shot (explodes), (very?) slow attack speed, high damage to enemies hit directly, lower damage to enemies further away from projectile when it collides with something)\" = \"Fireball (1 projectile per shot (explodes), (very?) slow attack speed, high damage to enemies hit directly, lower damage to enemies further away from projectile when it collides with something)\";\n}\n\n$tankOgreIndex = -1\n$burstFireIndex = -1\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($renames.ContainsKey($t)) {\n        $p.Range.Text = $renames[$t]\n        if ($t.StartsWith(\"Tank Enemy\")) {\n            $tankOgreIndex = $i\n        }\n    } elseif ($t.StartsWith(\"Burst-Fire (1 projectile per shot\")) {\n        $burstFireIndex = $i\n    }\n}\n\n# New enemy entries, inserted right after the (renamed) \"Staff Ogre\" entry.\nif ($tankOgreIndex -ge 0) {\n    $rng = $d.Paragraphs($tankOgreIndex).Range\n    $rng.InsertParagraphAfter()\n    $d.Paragraphs($tankOgreIndex + 1).Range.Text = \"Pyromaniac (Moves towards player quickly, uses Flamethrower, average health and fast movement)\"\n\n    $rng2 = $d.Paragraphs($tankOgreIndex + 1).Range\n    $rng2.InsertParagraphAfter()\n    $d.Paragraphs($tankOgreIndex + 2).Range.Text = \"Necromancer (Moves away from player slowly, spawns Standard Skeletons and shoots a ring of projectiles around itself occasionally, High health and slow movement)\"\n\n    # Two paragraphs were inserted before the weapon section, so shift the\n    # tracked Burst-Fire index down to keep it correct.\n    if ($burstFireIndex -gt $tankOgreIndex) {\n        $burstFireIndex = $burstFireIndex + 2\n    }\n}\n\n# New weapon entry, inserted right after the \"Burst-Fire\" entry.\nif ($burstFireIndex -ge 0) {\n    $rng3 = $d.Paragraphs($burstFireIndex).Range\n    $rng3.InsertParagraphAfter()\n    $d.Paragraphs($burstFireIndex + 1).Range.Text = \"Sword (Melee weapon, hits enemies directly in front of player, fast attack speed, high damage)\"\n}\n"}
